# Add data for 2022-03-21
# Updates the "through" date in the sheet name / header, and bumps a
# handful of neighborhood/month counts (plus a few brand-new data points)
# to reflect one additional day of data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet and update the header text to reflect the new "through" date.
$ws.Name = "Through 2022-03-13"
$ws.Range("B1").Value = "March 2022 (through March 13)"

# Bump existing counts by 1.
$ws.Range("K3").Value = 2    # Austin
$ws.Range("E4").Value = 4    # North Lawndale
$ws.Range("E5").Value = 3    # Garfield Park
$ws.Range("B6").Value = 2    # Rogers Park
$ws.Range("B13").Value = 2   # Woodlawn
$ws.Range("B15").Value = 2   # Humboldt Park
$ws.Range("B25").Value = 2   # Washington Park

# New data points that previously had no value.
$ws.Range("Q7").Value = 1    # South Shore
$ws.Range("Q9").Value = 1    # Chicago Lawn
$ws.Range("B47").Value = 1   # Brighton Park
$ws.Range("H57").Value = 1   # Douglas
$ws.Range("N75").Value = 1   # Oakland
$ws.Range("N84").Value = 1   # South Chicago
$ws.Range("H85").Value = 1   # South Deering
$ws.Range("B90").Value = 1   # Wrigleyville
